$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.35%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.46%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.700"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06101"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.02%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.707"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.63%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.12%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9134"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.82%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1405"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.51%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.83%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07079"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.15%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03127"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.78%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09053"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.02%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.70%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006165"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.41%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006128"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.54%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.449"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.10%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.171"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.61%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.166"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.14%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1300"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.79%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.133"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.08%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04255"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.45%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-3.65%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'6.44%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D40").Value = "'0.03952"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.48%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.05%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004179"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.38%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002134"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.78%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01314"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-19.44%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005125"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.58%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.02122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-61.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.2582"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'90.53%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E50").Style = "Normal"
